# Everyday data update: insert a new "today" row (2021/12/03) at the top of
# each of the 5 data sheets, pushing all existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
$ws1.Cells.Item(2,1).Value = "日期：2021/12/03"
$ws1.Cells.Item(2,2).NumberFormat = "@"
$ws1.Cells.Item(2,2).Value = "202201"
$ws1.Cells.Item(2,3).Value = 17650
$ws1.Cells.Item(2,4).Value = 9120
$ws1.Cells.Item(2,5).Value = 9478050
$ws1.Cells.Item(2,6).Value = 17633

# ---------------------------------------------------------------------
# Sheet 2: 散戶多空力道
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
$ws2.Cells.Item(2,1).Value = "日期：2021/12/03"
$ws2.Cells.Item(2,2).Value = 0.07000000000000001

# ---------------------------------------------------------------------
# Sheet 3: 三大法人買賣金額
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()
$ws3.Cells.Item(2,1).Value = "110年12月03日"
$ws3.Cells.Item(2,2).Value = 75.7
$ws3.Cells.Item(2,3).Value = -103.44

# ---------------------------------------------------------------------
# Sheet 4: 大盤多空點位
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
$ws4.Cells.Item(2,1).Value = "110年12月03日"
$ws4.Cells.Item(2,2).Value = 17728.92

# ---------------------------------------------------------------------
# Sheet 5: 期貨大額交易人未沖銷部位
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
$ws5.Cells.Item(2,1).NumberFormat = "@"
$ws5.Cells.Item(2,1).Value = "2021/12/03"
$ws5.Cells.Item(2,2).Value = 49287
$ws5.Cells.Item(2,3).Value = 53579
$ws5.Cells.Item(2,4).Value = 500
$ws5.Cells.Item(2,5).Value = 899
$ws5.Cells.Item(2,6).Value = 25658
$ws5.Cells.Item(2,7).Value = 47825
$ws5.Cells.Item(2,8).Value = -96
$ws5.Cells.Item(2,9).Value = 213
$ws5.Cells.Item(2,10).Value = -22167
$ws5.Cells.Item(2,11).Value = -309
$ws5.Cells.Item(2,12).Value = 596
$ws5.Cells.Item(2,13).Value = 686
$ws5.Cells.Item(2,14).Value = -90
